# Replace every occurrence of $old with $new, using Range.Text assignment
# (rather than Find.Execute's Replacement argument) so that Word's
# smart-quote AutoCorrect/AutoFormat does not mangle straight apostrophes.
function Replace-AllText($doc, $old, $new) {
    $count = 0
    $r = $doc.Content
    $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    while ($r.Find.Found) {
        $r.Text = $new
        $count = $count + 1
        $r.Collapse(0)
        $r.End = $doc.Content.End
        $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    }
    return $count
}

$d = $word.ActiveDocument

# Title (Heading1) and the bold repeated title near the end of the document
Replace-AllText $d "Play Gonzo's Quest Megaways Slot for Free - Review & Rating 2021" "Play Gonzo's Quest Megaways Free | Exciting Adventure Slot"

# "What we like" bullet list
Replace-AllText $d "Impressive graphics with mythical background" "Impressive graphics with stunning blurred effects"
Replace-AllText $d "Medium to high variance with 117,649 ways to win" "Exciting adventure theme in the ruins of El Dorado"
Replace-AllText $d "Multiple bonus features with avalanche reels, free fall bonus, and more" "Numerous bonus features and winning multipliers"
Replace-AllText $d "Max win potential of 21,000x" "High maximum win potential of 21,000x"

# "What we don't like" bullet list
Replace-AllText $d "Lowest bet amount is 10 cents" "Limited betting range with a maximum bet of €40.00 per spin"
Replace-AllText $d "The maximum bet amount is only €40.00 per spin" "Medium to high variance may not appeal to all players"

# Meta description (italic run near the end)
Replace-AllText $d "Read our in-depth review of Gonzo's Quest Megaways slot, and play for free. Check out our rating, bonus features, and max win potential of this Red Tiger Gaming slot." "Play Gonzo's Quest Megaways for free and experience an adventure in the ruins of El Dorado."
